# Auto-generated script applying scheduled-runner market data updates
# to the Adamantoise_Profits workbook (one hunk per Leve row, per job sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 15152227
$ws.Range("I92").Value = 17857870
$ws.Range("J92").Value = 628.4
$ws.Range("K92").Value = 17857870
$ws.Range("L92").Value = 628.4
$ws.Range("M92").Value = -17856622
$ws.Range("N92").Value = -3124.4

# Row 112
$ws.Range("H112").Value = 1496113.4
$ws.Range("I112").Value = 4300
$ws.Range("J112").Value = 1638190.9
$ws.Range("K112").Value = 12900
$ws.Range("L112").Value = 4914572.699999999
$ws.Range("M112").Value = -11792
$ws.Range("N112").Value = -4916788.699999999

# Row 132
$ws.Range("H132").Value = 2135.3914
$ws.Range("I132").Value = 2141.5908
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 6424.7724
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -3894.7724
$ws.Range("N132").Value = -11057

# Row 138
$ws.Range("H138").Value = 1775.77
$ws.Range("I138").Value = 850
$ws.Range("J138").Value = 2392.95
$ws.Range("K138").Value = 2550
$ws.Range("L138").Value = 7178.849999999999
$ws.Range("M138").Value = 2590
$ws.Range("N138").Value = -17458.85

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17965962
$ws.Range("I32").Value = 20489026
$ws.Range("J32").Value = 6497489.5
$ws.Range("K32").Value = 20489026
$ws.Range("L32").Value = 6497489.5
$ws.Range("M32").Value = -20488739
$ws.Range("N32").Value = -6498063.5

# Row 46
$ws.Range("H46").Value = 4025.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4025.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4025.5
$ws.Range("N46").Value = -4663.5

# Row 102
$ws.Range("H102").Value = 2087.25
$ws.Range("I102").Value = 1533
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 1533
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = 89
$ws.Range("N102").Value = -6994

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2826.5715
$ws.Range("I99").Value = 1939
$ws.Range("J99").Value = 3714.1428
$ws.Range("K99").Value = 1939
$ws.Range("L99").Value = 3714.1428
$ws.Range("M99").Value = -441
$ws.Range("N99").Value = -6710.1428

# Row 105
$ws.Range("H105").Value = 2387.3044
$ws.Range("I105").Value = 1911.0667
$ws.Range("J105").Value = 3280.25
$ws.Range("K105").Value = 1911.0667
$ws.Range("L105").Value = 3280.25
$ws.Range("M105").Value = -164.0667000000001
$ws.Range("N105").Value = -6774.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5256.6
$ws.Range("I31").Value = 2543.0476
$ws.Range("J31").Value = 7221.5864
$ws.Range("K31").Value = 2543.0476
$ws.Range("L31").Value = 7221.5864
$ws.Range("M31").Value = -2248.0476
$ws.Range("N31").Value = -7811.5864

# Row 34
$ws.Range("H34").Value = 5256.6
$ws.Range("I34").Value = 2543.0476
$ws.Range("J34").Value = 7221.5864
$ws.Range("K34").Value = 2543.0476
$ws.Range("L34").Value = 7221.5864
$ws.Range("M34").Value = -2341.0476
$ws.Range("N34").Value = -7625.5864

# Row 86
$ws.Range("H86").Value = 37030.426
$ws.Range("I86").Value = 37537.617
$ws.Range("J86").Value = 36700.75
$ws.Range("K86").Value = 37537.617
$ws.Range("L86").Value = 36700.75
$ws.Range("M86").Value = -36414.617
$ws.Range("N86").Value = -38946.75

# Row 89
$ws.Range("H89").Value = 37030.426
$ws.Range("I89").Value = 37537.617
$ws.Range("J89").Value = 36700.75
$ws.Range("K89").Value = 187688.085
$ws.Range("L89").Value = 183503.75
$ws.Range("M89").Value = -182072.085
$ws.Range("N89").Value = -194735.75

# Row 132
$ws.Range("H132").Value = 3425.3713
$ws.Range("I132").Value = 3046.2144
$ws.Range("J132").Value = 4942
$ws.Range("K132").Value = 9138.643199999999
$ws.Range("L132").Value = 14826
$ws.Range("M132").Value = -6608.643199999999
$ws.Range("N132").Value = -19886

# Row 134
$ws.Range("H134").Value = 2750
$ws.Range("I134").Value = 2666.6667
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 8000.000100000001
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -5465.000100000001
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 110158744
$ws.Range("I4").Value = 83743690
$ws.Range("J4").Value = 177663890
$ws.Range("K4").Value = 251231070
$ws.Range("L4").Value = 532991670
$ws.Range("M4").Value = -251230958
$ws.Range("N4").Value = -532991894

# Row 98
$ws.Range("H98").Value = 681.25
$ws.Range("I98").Value = 790.3333
$ws.Range("J98").Value = 615.8
$ws.Range("K98").Value = 2370.9999
$ws.Range("L98").Value = 1847.4
$ws.Range("M98").Value = -872.9998999999998
$ws.Range("N98").Value = -4843.4

$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 1250000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1250000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1250000
$ws.Range("N40").Value = -1250302

# Row 62
$ws.Range("H62").Value = 40000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 40000
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -41372

# Row 65
$ws.Range("H65").Value = 40000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 120000
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -126864

# Row 70
$ws.Range("H70").Value = 23587.53
$ws.Range("I70").Value = 67097.60000000001
$ws.Range("J70").Value = 5458.3335
$ws.Range("K70").Value = 67097.60000000001
$ws.Range("L70").Value = 5458.3335
$ws.Range("M70").Value = -66827.60000000001
$ws.Range("N70").Value = -5998.3335

# Row 73
$ws.Range("H73").Value = 23587.53
$ws.Range("I73").Value = 67097.60000000001
$ws.Range("J73").Value = 5458.3335
$ws.Range("K73").Value = 67097.60000000001
$ws.Range("L73").Value = 5458.3335
$ws.Range("M73").Value = -66161.60000000001
$ws.Range("N73").Value = -7330.3335

# Row 102
$ws.Range("H102").Value = 1435.9231
$ws.Range("I102").Value = 1354.3429
$ws.Range("J102").Value = 2149.75
$ws.Range("K102").Value = 1354.3429
$ws.Range("L102").Value = 2149.75
$ws.Range("M102").Value = 267.6570999999999
$ws.Range("N102").Value = -5393.75

# Row 132
$ws.Range("H132").Value = 10142.571
$ws.Range("I132").Value = 14333
$ws.Range("J132").Value = 6999.75
$ws.Range("K132").Value = 42999
$ws.Range("L132").Value = 20999.25
$ws.Range("M132").Value = -40469
$ws.Range("N132").Value = -26059.25

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2422.2222
$ws.Range("I22").Value = 1300.1818
$ws.Range("J22").Value = 4185.4287
$ws.Range("K22").Value = 1300.1818
$ws.Range("L22").Value = 4185.4287
$ws.Range("M22").Value = -1005.1818
$ws.Range("N22").Value = -4775.4287

# Row 27
$ws.Range("H27").Value = 2422.2222
$ws.Range("I27").Value = 1300.1818
$ws.Range("J27").Value = 4185.4287
$ws.Range("K27").Value = 1300.1818
$ws.Range("L27").Value = 4185.4287
$ws.Range("M27").Value = -1193.1818
$ws.Range("N27").Value = -4399.4287

# Row 55
$ws.Range("H55").Value = 263.46155
$ws.Range("I55").Value = 198
$ws.Range("J55").Value = 368.2
$ws.Range("K55").Value = 198
$ws.Range("L55").Value = 368.2
$ws.Range("M55").Value = -25
$ws.Range("N55").Value = -714.2

# Row 132
$ws.Range("H132").Value = 480869.84
$ws.Range("I132").Value = 628979.25
$ws.Range("J132").Value = 6919.8
$ws.Range("K132").Value = 1886937.75
$ws.Range("L132").Value = 20759.4
$ws.Range("M132").Value = -1884407.75
$ws.Range("N132").Value = -25819.4

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 30022.5
$ws.Range("I40").Value = 30025
$ws.Range("J40").Value = 30020
$ws.Range("K40").Value = 30025
$ws.Range("L40").Value = 30020
$ws.Range("M40").Value = -29876
$ws.Range("N40").Value = -30318

# Row 122
$ws.Range("H122").Value = 33337274
$ws.Range("I122").Value = 40004172
$ws.Range("J122").Value = 2791.6
$ws.Range("K122").Value = 120012516
$ws.Range("L122").Value = 8374.799999999999
$ws.Range("M122").Value = -120010066
$ws.Range("N122").Value = -13274.8

# Row 132
$ws.Range("H132").Value = 24988.342
$ws.Range("I132").Value = 29441.166
$ws.Range("J132").Value = 4950.625
$ws.Range("K132").Value = 88323.49800000001
$ws.Range("L132").Value = 14851.875
$ws.Range("M132").Value = -85793.49800000001
$ws.Range("N132").Value = -19911.875
